$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("U27").Value = 0.0308
$ws.Range("N42").Value = 0.3794
$ws.Range("O42").Value = 0.6053999999999999
$ws.Range("U42").Value = 0.0433
$ws.Range("V42").Value = 0.02
$ws.Range("W42").Value = 0.0011
$ws.Range("T43").Value = 0.0469
$ws.Range("N46").Value = 0.3892
$ws.Range("O46").Value = 0.6347
$ws.Range("U46").Value = 0.05309999999999999
$ws.Range("V46").Value = 0.02
$ws.Range("W46").Value = 0.0073
$ws.Range("N50").Value = 0.3746
$ws.Range("O50").Value = 0.6334000000000001
$ws.Range("U50").Value = 0.0731
$ws.Range("V50").Value = 0.02
$ws.Range("W50").Value = 0.001
$ws.Range("N54").Value = 0.361
$ws.Range("M58").Value = 0.4816
$ws.Range("N58").Value = 0.3515
$ws.Range("T58").Value = 0.0732
$ws.Range("M60").Value = 0.6663
$ws.Range("M62").Value = 0.4409
$ws.Range("N62").Value = 0.3290999999999999
$ws.Range("T62").Value = 0.0755
$ws.Range("V64").Value = 0.0401
$ws.Range("M66").Value = 0.4443
$ws.Range("N66").Value = 0.3356
$ws.Range("T66").Value = 0.0725
$ws.Range("V68").Value = 0.0401
$ws.Range("M70").Value = 0.4376
$ws.Range("V72").Value = 0.0401
